$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

$ws.Range("C7").Value = "sure bitti"

$ws.Range("G13").Select()
